$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 currently has the "empty last row" style (no value). Copy the
# formatting from the cell directly above it (B2) so B3 ends up sharing
# the same style slot as the rest of the data rows, then fill in
# Thabang's GitHub username.
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B3").Value = "thabanglwazi"
